$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H99").Value = 1494.5454
$ws.Range("I99").Value = 277.7143
$ws.Range("K99").Value = 833.1428999999999
$ws.Range("M99").Value = 664.8571000000001
$ws.Range("H103").Value = 1575
$ws.Range("J103").Value = 2500
$ws.Range("L103").Value = 7500
$ws.Range("N103").Value = -8672
$ws.Range("H106").Value = 1522.6364
$ws.Range("J106").Value = 1265
$ws.Range("L106").Value = 1265
$ws.Range("N106").Value = -2527
$ws.Range("H113").Value = 1500.9
$ws.Range("I113").Value = 1501
$ws.Range("K113").Value = 1501
$ws.Range("M113").Value = 1753
$ws.Range("H138").Value = 5533.5
$ws.Range("I138").Value = 4310.278
$ws.Range("J138").Value = 5831.0405
$ws.Range("K138").Value = 12930.834
$ws.Range("L138").Value = 17493.1215
$ws.Range("M138").Value = -7790.834000000001
$ws.Range("N138").Value = -27773.1215

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2505.5454
$ws.Range("I32").Value = 2505.5454
$ws.Range("K32").Value = 2505.5454
$ws.Range("M32").Value = -2218.5454
$ws.Range("H61").Value = 90911370
$ws.Range("I61").Value = 111112840
$ws.Range("J61").Value = 4750
$ws.Range("K61").Value = 111112840
$ws.Range("L61").Value = 4750
$ws.Range("M61").Value = -111112628
$ws.Range("N61").Value = -5174
$ws.Range("H63").Value = 3997.5
$ws.Range("I63").Value = 4163.3335
$ws.Range("J63").Value = 3500
$ws.Range("K63").Value = 4163.3335
$ws.Range("L63").Value = 3500
$ws.Range("M63").Value = -3477.3335
$ws.Range("N63").Value = -4872
$ws.Range("H66").Value = 3997.5
$ws.Range("I66").Value = 4163.3335
$ws.Range("J66").Value = 3500
$ws.Range("K66").Value = 20816.6675
$ws.Range("L66").Value = 17500
$ws.Range("M66").Value = -17384.6675
$ws.Range("N66").Value = -24364
$ws.Range("H74").Value = 111118140
$ws.Range("I74").Value = 166671380
$ws.Range("K74").Value = 166671380
$ws.Range("M74").Value = -166670506
$ws.Range("H77").Value = 111118140
$ws.Range("I77").Value = 166671380
$ws.Range("K77").Value = 833356900
$ws.Range("M77").Value = -833352532
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H136").Value = 90911370
$ws.Range("I136").Value = 111112840
$ws.Range("J136").Value = 4750
$ws.Range("K136").Value = 333338520
$ws.Range("L136").Value = 14250
$ws.Range("M136").Value = -333335970
$ws.Range("N136").Value = -19350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2165.3333
$ws.Range("I99").Value = 1999.5
$ws.Range("J99").Value = 2248.25
$ws.Range("K99").Value = 1999.5
$ws.Range("L99").Value = 2248.25
$ws.Range("M99").Value = -501.5
$ws.Range("N99").Value = -5244.25
$ws.Range("H107").Value = 335146
$ws.Range("I107").Value = 2439
$ws.Range("J107").Value = 501499.5
$ws.Range("K107").Value = 2439
$ws.Range("L107").Value = 501499.5
$ws.Range("M107").Value = -519
$ws.Range("N107").Value = -505339.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2507.25
$ws.Range("I16").Value = 1722.5714
$ws.Range("K16").Value = 1722.5714
$ws.Range("M16").Value = -1435.5714
$ws.Range("H86").Value = 7775.091
$ws.Range("I86").Value = 6602.6
$ws.Range("K86").Value = 6602.6
$ws.Range("M86").Value = -5479.6
$ws.Range("H89").Value = 7775.091
$ws.Range("I89").Value = 6602.6
$ws.Range("K89").Value = 33013
$ws.Range("M89").Value = -27397
$ws.Range("H99").Value = 2575.4
$ws.Range("I99").Value = 2469.25
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 2469.25
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -971.25
$ws.Range("N99").Value = -5996
$ws.Range("H107").Value = 77563.30499999999
$ws.Range("I107").Value = 693.0909
$ws.Range("K107").Value = 693.0909
$ws.Range("M107").Value = 1226.9091
$ws.Range("H113").Value = 2507.25
$ws.Range("I113").Value = 1722.5714
$ws.Range("K113").Value = 1722.5714
$ws.Range("M113").Value = 447.4286
$ws.Range("H122").Value = 3699.8235
$ws.Range("I122").Value = 3671.2856
$ws.Range("K122").Value = 11013.8568
$ws.Range("M122").Value = -8563.856800000001
$ws.Range("H126").Value = 2575.4
$ws.Range("I126").Value = 2469.25
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 7407.75
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -4937.75
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 984.9
$ws.Range("I97").Value = 606
$ws.Range("J97").Value = 1111.2
$ws.Range("K97").Value = 1818
$ws.Range("L97").Value = 3333.6
$ws.Range("M97").Value = -1322
$ws.Range("N97").Value = -4325.6
$ws.Range("H107").Value = 1082.8572
$ws.Range("I107").Value = 790
$ws.Range("J107").Value = 1174.375
$ws.Range("K107").Value = 2370
$ws.Range("L107").Value = 3523.125
$ws.Range("M107").Value = -450
$ws.Range("N107").Value = -7363.125
$ws.Range("H131").Value = 2514.7778
$ws.Range("J131").Value = 2529.125
$ws.Range("L131").Value = 7587.375
$ws.Range("N131").Value = -17667.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2943
$ws.Range("I80").Value = 3089.6
$ws.Range("J80").Value = 2698.6667
$ws.Range("K80").Value = 3089.6
$ws.Range("L80").Value = 2698.6667
$ws.Range("M80").Value = -2091.6
$ws.Range("N80").Value = -4694.6667
$ws.Range("H83").Value = 2943
$ws.Range("I83").Value = 3089.6
$ws.Range("J83").Value = 2698.6667
$ws.Range("K83").Value = 15448
$ws.Range("L83").Value = 13493.3335
$ws.Range("M83").Value = -10456
$ws.Range("N83").Value = -23477.3335
$ws.Range("H102").Value = 4696.5557
$ws.Range("I102").Value = 2023.5454
$ws.Range("K102").Value = 2023.5454
$ws.Range("M102").Value = -401.5454
$ws.Range("H132").Value = 2372927.5
$ws.Range("I132").Value = 2417599.2
$ws.Range("K132").Value = 7252797.600000001
$ws.Range("M132").Value = -7250267.600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2833.1667
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H34").Value = 5000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H40").Value = 4079.4
$ws.Range("I40").Value = 4079.4
$ws.Range("K40").Value = 4079.4
$ws.Range("M40").Value = -3943.4
$ws.Range("H93").Value = 1210.3125
$ws.Range("I93").Value = 1211
$ws.Range("K93").Value = 1211
$ws.Range("M93").Value = 37
$ws.Range("H106").Value = 23500
$ws.Range("J106").Value = 23500
$ws.Range("L106").Value = 23500
$ws.Range("N106").Value = -26024
$ws.Range("H122").Value = 4498.1665
$ws.Range("I122").Value = 4598
$ws.Range("K122").Value = 13794
$ws.Range("M122").Value = -11344
$ws.Range("H126").Value = 2833.1667
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1249.5
$ws.Range("J81").Value = 1000
$ws.Range("L81").Value = 2000
$ws.Range("N81").Value = -4122
$ws.Range("H84").Value = 1249.5
$ws.Range("J84").Value = 1000
$ws.Range("L84").Value = 10000
$ws.Range("N84").Value = -20608
$ws.Range("H104").Value = 39000
$ws.Range("J104").Value = 39000
$ws.Range("L104").Value = 39000
$ws.Range("N104").Value = -45988
$ws.Range("H122").Value = 3120.6667
$ws.Range("J122").Value = 2299
$ws.Range("L122").Value = 6897
$ws.Range("N122").Value = -11797
$ws.Range("H127").Value = 90195
$ws.Range("J127").Value = 100000
$ws.Range("L127").Value = 100000
$ws.Range("N127").Value = -109920

